# Generate Report for Handoff
# Update status text from "In Translation" to "Ready for handoff" and
# refresh the associated handoff timestamps across the Overview, zh-cn,
# and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-16 18:53:19"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-16 18:53:14"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-16 18:53:19"

# Widen the status columns now that the text is longer, matching the
# column width growth recorded for the edit (~13.41 -> ~17.22 "characters").
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
